$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.391.34"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "3.062.47"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'592.75"
$ws.Range("D6").Value = "'154.23"
$ws.Range("E6").Value = "  +1.16%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -3.56%  "
$ws.Range("D9").Value = "3.061.86"
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("D10").Value = "'0.155"
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("E12").Value = "  -2.42%  "
$ws.Range("D13").Value = "'36.83"
$ws.Range("E13").Value = "  -1.32%  "
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "3.572.83"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").Value = "63.367.42"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").Value = "3.068.44"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D20").Value = "'489.17"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").Value = "'14.39"
$ws.Range("E21").Value = "  -1.99%  "
$ws.Range("E22").Value = "  -1.70%  "
$ws.Range("E24").Value = "  +4.02%  "
$ws.Range("D25").Value = "'82.12"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("D26").Value = "'12.89"
$ws.Range("E26").Value = "  -1.67%  "
$ws.Range("D27").Value = "'10.72"
$ws.Range("E27").Value = "  +10.65%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "'7.38"
$ws.Range("E29").Value = "  +2.05%  "
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").Value = "'27.43"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("E34").Value = "  -1.29%  "
$ws.Range("D35").Value = "'1.07"
$ws.Range("E35").Value = "  +1.01%  "
$ws.Range("D36").Value = "0.0₃0823"
$ws.Range("E36").Value = "  -3.68%  "
$ws.Range("D37").Value = "'3.34"
$ws.Range("E37").Value = "  -0.59%  "
$ws.Range("D38").Value = "'5.99"
$ws.Range("E38").Value = "  -2.41%  "
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("D41").Value = "'50.61"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("D42").Value = "'439.24"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").Value = "'0.292"
$ws.Range("E43").Value = "  +2.46%  "
$ws.Range("E44").Value = "  +3.03%  "
$ws.Range("D45").Value = "'0.0364"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "2.851.88"
$ws.Range("E46").Value = "  +1.30%  "
$ws.Range("D47").Value = "'38.89"
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("D48").Value = "'130.28"
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'25.29"
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("E51").Value = "  -1.04%  "
